# MifosLab Installation.docx - "Updates to lab environment"
#
# 1) "helm del --purge "            -> "helm uninstall "
# 2) "helm search -l "              -> "helm search repo "
# 3) "helm --namespace demo --name moja install mojaloop/mojaloop --version"
#       -> "helm install moja mojaloop/mojaloop--namespace demo --version"
#    (the --namespace/--name flags move from the front of the command to a
#     new "--namespace demo" token tacked onto the end, right before
#     " --version <version>"; this is also where the _GoBack bookmark,
#     which Word leaves at the point of the most recent edit, now lives)

$d = $word.ActiveDocument

# --- 1) helm del --purge -> helm uninstall --------------------------------
$r1 = $d.Content
$ok1 = $r1.Find.Execute("helm del --purge ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "helm uninstall ", 2)

# --- 2) helm search -l -> helm search repo --------------------------------
$r2 = $d.Content
$ok2 = $r2.Find.Execute("helm search -l ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "helm search repo ", 2)

# --- 3) helm --namespace demo --name  -> helm install ---------------------
$r3 = $d.Content
$ok3 = $r3.Find.Execute("helm --namespace demo --name ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "helm install ", 2)

# --- 4) drop the old " install " between "moja" and "mojaloop" ------------
$r4 = $d.Content
$ok4 = $r4.Find.Execute("moja install mojaloop", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "moja mojaloop", 2)

# --- 5) append the relocated --namespace flag after the final "mojaloop" --
$r5 = $d.Content
$ok5 = $r5.Find.Execute("mojaloop/mojaloop --version", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "mojaloop/mojaloop--namespace demo --version", 2)

# --- 6) relocate the _GoBack bookmark to the new edit point ---------------
# It currently sits on its own empty paragraph near the end of the document
# (left over from whatever paragraph Word considers "last edited"); move it
# to just after the newly-inserted "--namespace demo" text.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$anchor = $d.Content
$anchor.Find.Execute("--namespace demo --version", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$bmStart = $anchor.Start + 16   # length of "--namespace demo"
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "helm del->uninstall: $ok1"
Write-Host "helm search -l->repo: $ok2"
Write-Host "helm namespace/name->install: $ok3"
Write-Host "drop inner install: $ok4"
Write-Host "append namespace flag: $ok5"
